$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.230.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4691"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07937"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.169"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6807"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.221.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("E18").Value = "  +8.59%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007384"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.113.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.319"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.07%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.184"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.215"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.959"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09842"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.372"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.051"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04705"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7036"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.708"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01870"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.613"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.293"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.936"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8457"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4162"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "955.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.153"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.169"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3899"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.17%  "
